# Updated symbol list on Wed Dec 28 02:38:01 UTC 2022 with GitHub Actions
#
# Price (column D) values are stored as text in the workbook even though
# they look numeric, so a leading apostrophe is used to force Excel to
# keep them as text (matching the original formatting) instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.92"
$ws.Range("D3").Value = "'24.04"
$ws.Range("D4").Value = "'5.354"
$ws.Range("D5").Value = "'0.05799"
$ws.Range("D6").Value = "'6.473"
$ws.Range("D7").Value = "'3.337"
$ws.Range("D8").Value = "'0.8086"
$ws.Range("D9").Value = "'0.9216"
$ws.Range("D10").Value = "'0.1395"
$ws.Range("D11").Value = "'0.07393"
$ws.Range("D12").Value = "'0.03175"
$ws.Range("D13").Value = "'0.03059"
$ws.Range("D14").Value = "'0.09356"
$ws.Range("D15").Value = "'3.868"
$ws.Range("D16").Value = "'0.001559"
$ws.Range("D17").Value = "'0.04754"
$ws.Range("D18").Value = "'0.0005994"
$ws.Range("D19").Value = "'0.005900"
$ws.Range("D21").Value = "'0.004654"
$ws.Range("D22").Value = "'0.00008801"
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$ws.Range("D23").Value = "'3.613"
$ws.Range("D24").Value = "'2.141"
$ws.Range("D25").Value = "'0.3178"
$ws.Range("E28").Value = '27UpBotsUBXT'
$ws.Range("D40").Value = "'0.03831"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.006371"
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1065"
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.002750"
$ws.Range("E43").Value = '42CEJICEJIWorstin24h'
$ws.Range("D44").Value = "'0.008271"
$ws.Range("D45").Value = "'0.00005320"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D48").Value = "'0.001843"
